$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The log sheet already has entries in rows 2..72; append the newest run
# result as the next row, copying the formatting from the prior row.
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
$newRow = $lastRow + 1

$srcRange = $ws.Range($ws.Cells.Item($lastRow, 1), $ws.Cells.Item($lastRow, 8))
$dstRange = $ws.Range($ws.Cells.Item($newRow, 1), $ws.Cells.Item($newRow, 8))
$srcRange.Copy()
$dstRange.PasteSpecial(-4122)

$ws.Cells.Item($newRow, 1).Value = "2025-08-29 13:01:38 UTC"
$ws.Cells.Item($newRow, 2).Value = "2025-08-29 18:31:38 IST"
$ws.Cells.Item($newRow, 3).Value = "SKIPPED"
$ws.Cells.Item($newRow, 4).Value = "No change in PDF. Skipping download & Excel update."
$ws.Cells.Item($newRow, 5).Value = "https://nalcoindia.com/wp-content/uploads/2019/01/INGOT-28-08-2025.pdf"
$ws.Cells.Item($newRow, 6).Value = ""
$ws.Cells.Item($newRow, 7).Value = 0
$ws.Cells.Item($newRow, 8).Value = ""
